$wb = $excel.ActiveWorkbook

# --- Rename first sheet, insert new "Lookup" sheet between it and "Data" ---
$ws1 = $wb.Worksheets.Item("PoiFormulaHelperTest")
$ws1.Name = "ExcelFormulaTreeTest"

$wsLookup = $wb.Worksheets.Add($wb.Worksheets.Item("Data"))
$wsLookup.Name = "Lookup"

# NOTE: worksheet handles in this host resolve by index, and inserting a
# sheet shifts indices, so re-fetch "Data" (and even "ExcelFormulaTreeTest")
# by name now that the sheet order has settled.
$ws1 = $wb.Worksheets.Item("ExcelFormulaTreeTest")
$wsData = $wb.Worksheets.Item("Data")

# ======================================================================
# Sheet "ExcelFormulaTreeTest" - add label column + more formula rows
# (fill order chosen to reproduce the shared-string table ordering)
# ======================================================================
$ws1.Range("B3").Value = "Arithmetic no brackets"
$ws1.Range("B4").Value = "Arithmetic 1 bracket"
$ws1.Range("B5").Value = "Arithmetic 1 bracket variation"
$ws1.Range("B6").Value = "Aritmetic 2 sets of brackets"
$ws1.Range("B1").Value = "SUM over 1 operand or"
$ws1.Range("B2").Value = "SUM over multiple operands"
$ws1.Range("B7").Value = "Unary operation formula"
$ws1.Range("B8").Value = "Percentage formula"
$ws1.Range("B9").Value = "Multiple Function Eval"
$ws1.Range("B10").Value = "VLOOKUP"

$ws1.Range("A3").Formula = "=65+20"
$ws1.Range("A4").Formula = "=(6*5)+500"
$ws1.Range("A5").Formula = "=6*(5+500)"
$ws1.Range("A6").Formula = "=(34*45)+(800/40)"
$ws1.Range("A7").Formula = "=-(8-4)"
$ws1.Range("A8").Formula = "=9%"
$ws1.Range("A9").Formula = "=SUM(Data!A1:D5)+SUM(Data!B1:D5)"
$ws1.Range("A10").Formula = '=VLOOKUP("A1", Lookup!A1:E5, 3)'
$ws1.Range("A11").Formula = '=VLOOKUP("A2", Lookup!A1:E5, 3)'

# ======================================================================
# Sheet "Lookup" - lookup table (fill order reproduces shared-string order)
# ======================================================================
$wsLookup.Range("A1").Value = "A1"
$wsLookup.Range("A2").Value = "A2"
$wsLookup.Range("A3").Value = "A3"
$wsLookup.Range("A4").Value = "A4"
$wsLookup.Range("A5").Value = "A5"

$wsLookup.Range("C1").Value = "C1 Result"
$wsLookup.Range("D1").Value = "D1 Result"
$wsLookup.Range("E1").Value = "E1 Result"
$wsLookup.Range("B1").Value = "B1 Result"

$wsLookup.Range("B2").Value = "B2 Result"
$wsLookup.Range("C2").Value = "C2 Result"
$wsLookup.Range("D2").Value = "D2 Result"
$wsLookup.Range("E2").Value = "E2 Result"

$wsLookup.Range("B3").Value = "B3 Result"
$wsLookup.Range("C3").Value = "C3 Result"
$wsLookup.Range("D3").Value = "D3 Result"
$wsLookup.Range("E3").Value = "E3 Result"

$wsLookup.Range("B4").Value = "B4 Result"
$wsLookup.Range("C4").Value = "C4 Result"
$wsLookup.Range("D4").Value = "D4 Result"
$wsLookup.Range("E4").Value = "E4 Result"

$wsLookup.Range("B5").Value = "B5 Result"
$wsLookup.Range("C5").Value = "C5 Result"
$wsLookup.Range("D5").Value = "D5 Result"
$wsLookup.Range("E5").Value = "E5 Result"

# --- Selections matching the final saved view state ---
$wsLookup.Select() | Out-Null
$wsLookup.Range("F17").Select() | Out-Null

$wsData.Select() | Out-Null
$wsData.Range("A1:D5").Select() | Out-Null

$ws1.Select() | Out-Null
$ws1.Range("A9").Select() | Out-Null
